# Fruta / hortaliza, semanal
# Reassign the Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M), Origen (O) and Precio $/Kg (P) values across
# the data rows (2-14) as per the updated weekly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = 44189; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota";   P = 344 },
    @{ Row = 3;  D = 44251; J = 120; K = 5000; L = 5000; M = 5000; O = "Región Metropolitana";     P = 312 },
    @{ Row = 4;  D = 44208; J = 160; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota";   P = 344 },
    @{ Row = 5;  D = 44210; J = 340; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota";   P = 344 },
    @{ Row = 6;  D = 44188; J = 210; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota";   P = 344 },
    @{ Row = 7;  D = 44230; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota";   P = 344 },
    @{ Row = 8;  D = 44232; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota";   P = 344 },
    @{ Row = 9;  D = 44215; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota";   P = 344 },
    @{ Row = 10; D = 44292; J = 90;  K = 6000; L = 6000; M = 6000; O = "Región Metropolitana";     P = 375 },
    @{ Row = 11; D = 44231; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota";   P = 344 },
    @{ Row = 12; D = 44204; J = 430; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota";   P = 344 },
    @{ Row = 13; D = 44186; J = 160; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota";   P = 344 },
    @{ Row = 14; D = 44187; J = 160; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota";   P = 344 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value = $u.D    # D: Fecha
    $ws.Cells.Item($r, 10).Value = $u.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $u.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $u.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $u.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value = $u.O   # O: Origen
    $ws.Cells.Item($r, 16).Value = $u.P   # P: Precio $/Kg
}
